$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$co = $ws.ChartObjects().Add(100, 100, 300, 200)
$chart = $co.Chart
$chart.ChartType = 5
try {
  $shp = $chart.Shapes.AddTextbox(1, 10, 10, 100, 30)
  $shp.TextFrame.Characters().Text = "sin detecciones inválidas"
  Write-Output "added textbox"
} catch {
  Write-Output "ERR: $_"
}
